$d = $word.ActiveDocument

# 1) Split the "Version" run into "Versi" + "on" without introducing any
#    direct character formatting: temporarily drop a bookmark at the split
#    point (char offset 5 -> right after "Versi"), which forces the host
#    to break the run there, then delete the bookmark again.
$splitPoint = $d.Range(5, 5)
$d.Bookmarks.Add("_tmpSplit", $splitPoint) | Out-Null
$d.Bookmarks("_tmpSplit").Delete()

# 2) "Version 1." -> "Version 2" (keep the trailing "." out of this run;
#    it gets reinserted after the _GoBack bookmark in step 3).
$d.Content.Find.Execute(" 1.", $true, $false, $false, $false, $false, $true, 1, $false, " 2", 2)

# 3) Re-append the final "." as its own run after the _GoBack bookmark.
$endOfPara = $d.Content.Text.Length - 1
$tail = $d.Range($endOfPara, $endOfPara)
$tail.InsertAfter(".")
